# Generate Report for Handback
#
# The handback for 8f104349-c4d2-4df8-be52-d8076a42e2d6 finished and is now
# in sync with en-US, so the status report is regenerated:
#   - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview!E2/F2, zh-cn!C2, de-de!C2)
#   - The handback timestamps for both locales move forward
#   - The stale "handback file is not latest" error is cleared now that the
#     handback is current
#   - A handful of columns get re-measured to fit the new text

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# Columns widened to fit the longer status text.
$overview.Range("E1").ColumnWidth = 29.15
$overview.Range("F1").ColumnWidth = 29.15

# --- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("K2").Value = "2016-09-04 12:53:30"
$zhcn.Range("P2").Value = ""

$zhcn.Range("C1").ColumnWidth = 29.15
$zhcn.Range("P1").ColumnWidth = 12.85

# --- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("K2").Value = "2016-09-04 12:53:38"
$dede.Range("P2").Value = ""

$dede.Range("C1").ColumnWidth = 29.15
$dede.Range("P1").ColumnWidth = 12.85
